$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1 so the new header cell matches the existing header formatting
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

$timeTaken = @(
  "2021-10-05 13:40:16.132137",
  "2021-10-05 13:40:16.132150",
  "2021-10-05 13:40:16.132154",
  "2021-10-05 13:40:16.132157",
  "2021-10-05 13:40:16.132161",
  "2021-10-05 13:40:16.132164",
  "2021-10-05 13:40:16.132167",
  "2021-10-05 13:40:16.132170",
  "2021-10-05 13:40:16.132174",
  "2021-10-05 13:40:16.132177",
  "2021-10-05 13:40:16.132180",
  "2021-10-05 13:40:16.132183",
  "2021-10-05 13:40:16.132186",
  "2021-10-05 13:40:16.132189",
  "2021-10-05 13:40:16.132193",
  "2021-10-05 13:40:16.132196",
  "2021-10-05 13:40:16.132199",
  "2021-10-05 13:40:16.132203",
  "2021-10-05 13:40:16.132206",
  "2021-10-05 13:40:16.132209",
  "2021-10-05 13:40:16.132212",
  "2021-10-05 13:40:16.132215",
  "2021-10-05 13:40:16.132218",
  "2021-10-05 13:40:16.132221",
  "2021-10-05 13:40:16.132225",
  "2021-10-05 13:40:16.132228",
  "2021-10-05 13:40:16.132231",
  "2021-10-05 13:40:16.132234",
  "2021-10-05 13:40:16.132238",
  "2021-10-05 13:40:16.132241",
  "2021-10-05 13:40:16.132244",
  "2021-10-05 13:40:16.132247",
  "2021-10-05 13:40:16.132251",
  "2021-10-05 13:40:16.132254",
  "2021-10-05 13:40:16.132257",
  "2021-10-05 13:40:16.132261",
  "2021-10-05 13:40:16.132264",
  "2021-10-05 13:40:16.132267",
  "2021-10-05 13:40:16.132270",
  "2021-10-05 13:40:16.132273",
  "2021-10-05 13:40:16.132277",
  "2021-10-05 13:40:16.132281",
  "2021-10-05 13:40:16.132284",
  "2021-10-05 13:40:16.132287",
  "2021-10-05 13:40:16.132290",
  "2021-10-05 13:40:16.132294",
  "2021-10-05 13:40:16.132297",
  "2021-10-05 13:40:16.132300",
  "2021-10-05 13:40:16.132303",
  "2021-10-05 13:40:16.132306",
  "2021-10-05 13:40:16.132309",
  "2021-10-05 13:40:16.132313",
  "2021-10-05 13:40:16.132316",
  "2021-10-05 13:40:16.132320",
  "2021-10-05 13:40:16.132323",
  "2021-10-05 13:40:16.132326",
  "2021-10-05 13:40:16.132329",
  "2021-10-05 13:40:16.132333",
  "2021-10-05 13:40:16.132336",
  "2021-10-05 13:40:16.132339",
  "2021-10-05 13:40:16.132342",
  "2021-10-05 13:40:16.132345",
  "2021-10-05 13:40:16.132348",
  "2021-10-05 13:40:16.132351",
  "2021-10-05 13:40:16.132355",
  "2021-10-05 13:40:16.132359",
  "2021-10-05 13:40:16.132362",
  "2021-10-05 13:40:16.132365",
  "2021-10-05 13:40:16.132368",
  "2021-10-05 13:40:16.132371",
  "2021-10-05 13:40:16.132374",
  "2021-10-05 13:40:16.132377",
  "2021-10-05 13:40:16.132380",
  "2021-10-05 13:40:16.132383",
  "2021-10-05 13:40:16.132386",
  "2021-10-05 13:40:16.132389",
  "2021-10-05 13:40:16.132394",
  "2021-10-05 13:40:16.132397",
  "2021-10-05 13:40:16.132400",
  "2021-10-05 13:40:16.132403",
  "2021-10-05 13:40:16.132407",
  "2021-10-05 13:40:16.132410",
  "2021-10-05 13:40:16.132412",
  "2021-10-05 13:40:16.132415",
  "2021-10-05 13:40:16.132418",
  "2021-10-05 13:40:16.132421",
  "2021-10-05 13:40:16.132424",
  "2021-10-05 13:40:16.132428",
  "2021-10-05 13:40:16.132430",
  "2021-10-05 13:40:16.132433",
  "2021-10-05 13:40:16.132436",
  "2021-10-05 13:40:16.132439",
  "2021-10-05 13:40:16.132443",
  "2021-10-05 13:40:16.132447",
  "2021-10-05 13:40:16.132450",
  "2021-10-05 13:40:16.132453",
  "2021-10-05 13:40:16.132456",
  "2021-10-05 13:40:16.132459",
  "2021-10-05 13:40:16.132462",
  "2021-10-05 13:40:16.132465",
  "2021-10-05 13:40:16.132468",
  "2021-10-05 13:40:16.132471",
  "2021-10-05 13:40:16.132474",
  "2021-10-05 13:40:16.132477",
  "2021-10-05 13:40:16.132480",
  "2021-10-05 13:40:16.132483",
  "2021-10-05 13:40:16.132486",
  "2021-10-05 13:40:16.132489",
  "2021-10-05 13:40:16.132494",
  "2021-10-05 13:40:16.132498",
  "2021-10-05 13:40:16.132501",
  "2021-10-05 13:40:16.132504",
  "2021-10-05 13:40:16.132507",
  "2021-10-05 13:40:16.132510",
  "2021-10-05 13:40:16.132513",
  "2021-10-05 13:40:16.132516",
  "2021-10-05 13:40:16.132519",
  "2021-10-05 13:40:16.132522",
  "2021-10-05 13:40:16.132525",
  "2021-10-05 13:40:16.132528",
  "2021-10-05 13:40:16.132531",
  "2021-10-05 13:40:16.132534",
  "2021-10-05 13:40:16.132537"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

